$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Branch (column D) for existing rows per diff
$ws.Range("D8").Value = "CS"
$ws.Range("D9").Value = "CS"
$ws.Range("D11").Value = "IT"
$ws.Range("D13").Value = "IT"

# Add new row 14
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 85
$ws.Range("C14").Value = "Pratik"
$ws.Range("D14").Value = "CS"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = "B"
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = "+"

# Match style of B column (centered) for new row's roll-no cell
$ws.Range("B14").HorizontalAlignment = -4108

# Update selection to reflect new active cell
$ws.Range("H14").Select()
